$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "is_checked" column (G) — header + "no" for every data row.
$ws.Range("G1").Value = "is_checked"
for ($r = 2; $r -le 42; $r++) {
    $ws.Range("G$r").Value = "no"
}

# Update the view's active cell / selection to match the edited workbook.
$ws.Range("I33").Select() | Out-Null
